$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so that
# numeric-looking values (e.g. "1.004") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.614.46'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.874.32'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '325.09'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4542'
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("D8").Value = '0.3849'
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '46.53'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07838'
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '0.9922'
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '21.54'
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.897.64'
$ws.Range("E13").Value = '  +4.00%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '6.931'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.647'
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = '0.06933'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '86.75'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = '1.006'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000009962'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '16.74'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.611.46'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.266'
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.069'
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.098.79'
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '154.42'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '19.15'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '5.683'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '117.46'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("B31").Value = 'LidoDAOToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D31").Value = '1.890'
$ws.Range("E31").Value = '  -5.13%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.09260'
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '0.9065'
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.282'
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.324'
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '3.256'
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.05674'
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.149'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02042'
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '7.635'
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5577'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1771'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '9.654'
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '0.07159'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").Value = '0.5261'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").Value = '2.157'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '11.48'
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '1.116'
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.809'
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '111.76'
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '2.438'
$ws.Range("E51").Value = '  +4.07%  '
